# Update "想去人数" (interest count) values in column F across the
# "展览" (sheet1), "本地生活" (sheet3) and "全部类型" (sheet4) sheets.
# "演出" (sheet2) is untouched.

$wb = $excel.ActiveWorkbook

# --- 展览 ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value  = 5312
$ws1.Range("F7").Value  = 456
$ws1.Range("F8").Value  = 643
$ws1.Range("F9").Value  = 920
$ws1.Range("F13").Value = 573
$ws1.Range("F17").Value = 1778
$ws1.Range("F18").Value = 1458
$ws1.Range("F19").Value = 841
$ws1.Range("F22").Value = 314
$ws1.Range("F24").Value = 138
$ws1.Range("F28").Value = 2638
$ws1.Range("F34").Value = 291
$ws1.Range("F35").Value = 11
$ws1.Range("F39").Value = 277
$ws1.Range("F40").Value = 649

# --- 本地生活 ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 235

# --- 全部类型 ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 235
$ws4.Range("F7").Value  = 5312
$ws4.Range("F8").Value  = 456
$ws4.Range("F9").Value  = 643
$ws4.Range("F12").Value = 920
$ws4.Range("F18").Value = 573
$ws4.Range("F23").Value = 1778
$ws4.Range("F24").Value = 1458
$ws4.Range("F25").Value = 841
$ws4.Range("F27").Value = 314
$ws4.Range("F30").Value = 138
$ws4.Range("F32").Value = 2638
$ws4.Range("F38").Value = 291
$ws4.Range("F39").Value = 11
$ws4.Range("F42").Value = 277
$ws4.Range("F43").Value = 649
